$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Sending cluster (column A) for the two row-groups
$ws.Range("A2:A6").Value = "FAPs"
$ws.Range("A7:A11").Value = "Resolving-Mac"

# Row 2
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.4870623333333333
$ws.Range("H2").Value = 1.461187
$ws.Range("I2").Value = 0.8073404988294784
$ws.Range("J2").Value = 0.8073404988294784
$ws.Range("M2").Value = 0.03269333333333333
$ws.Range("N2").Value = 0.09808
$ws.Range("O2").Value = 0.001656893562427925
$ws.Range("P2").Value = 0.001659236418317591
$ws.Range("Q2").Value = 0.01592369121777778
$ws.Range("R2").Value = 0.14331322096
$ws.Range("S2").Value = 0.001337677275197912
$ws.Range("T2").Value = 0.001339568757640561

# Row 3
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.4870623333333333
$ws.Range("H3").Value = 1.461187
$ws.Range("I3").Value = 0.8073404988294784
$ws.Range("J3").Value = 0.8073404988294784
$ws.Range("O3").Value = 0.09959606124235393
$ws.Range("P3").Value = 0.09973689057741834
$ws.Range("Q3").Value = 0.9571748974666666
$ws.Range("R3").Value = 8.6145740772
$ws.Range("S3").Value = 0.08040793376485331
$ws.Range("T3").Value = 0.08052163099047403

# Row 4
$ws.Range("D4").Value = "Inflammatory-Mac"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.4870623333333333
$ws.Range("H4").Value = 1.461187
$ws.Range("I4").Value = 0.8073404988294784
$ws.Range("J4").Value = 0.8073404988294784
$ws.Range("M4").Value = 13.379326
$ws.Range("N4").Value = 40.137978
$ws.Range("O4").Value = 0.6780623710957755
$ws.Range("P4").Value = 0.679021154722984
$ws.Range("Q4").Value = 6.516565739987333
$ws.Range("R4").Value = 58.649091659886
$ws.Range("S4").Value = 0.5474272129179624
$ws.Range("T4").Value = 0.5482012777698223

# Row 5
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.4870623333333333
$ws.Range("H5").Value = 1.461187
$ws.Range("I5").Value = 0.8073404988294784
$ws.Range("J5").Value = 0.8073404988294784
$ws.Range("M5").Value = 0.08358400000000001
$ws.Range("N5").Value = 0.167168
$ws.Range("O5").Value = 0.004236025433991914
$ws.Range("P5").Value = 0.00282801013027442
$ws.Range("Q5").Value = 0.04071061806933333
$ws.Range("R5").Value = 0.244263708416
$ws.Range("S5").Value = 0.00341991488693339
$ws.Range("T5").Value = 0.002283167109270569

# Row 6
$ws.Range("D6").Value = "Resolving-Mac"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.4870623333333333
$ws.Range("H6").Value = 1.461187
$ws.Range("I6").Value = 0.8073404988294784
$ws.Range("J6").Value = 0.8073404988294784
$ws.Range("M6").Value = 4.270900666666667
$ws.Range("N6").Value = 12.812702
$ws.Range("O6").Value = 0.2164486486654506
$ws.Range("P6").Value = 0.2167547081510057
$ws.Range("Q6").Value = 2.080194844141555
$ws.Range("R6").Value = 18.721753597274
$ws.Range("S6").Value = 0.1747477599845314
$ws.Range("T6").Value = 0.174994854202271

# Row 7
$ws.Range("D7").Value = "ECs"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.11623
$ws.Range("H7").Value = 0.34869
$ws.Range("I7").Value = 0.1926595011705215
$ws.Range("J7").Value = 0.1926595011705215
$ws.Range("M7").Value = 0.03269333333333333
$ws.Range("N7").Value = 0.09808
$ws.Range("O7").Value = 0.001656893562427925
$ws.Range("P7").Value = 0.001659236418317591
$ws.Range("Q7").Value = 0.003799946133333333
$ws.Range("R7").Value = 0.0341995152
$ws.Range("S7").Value = 0.0003192162872300123
$ws.Range("T7").Value = 0.0003196676606770299

# Row 8
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.11623
$ws.Range("H8").Value = 0.34869
$ws.Range("I8").Value = 0.1926595011705215
$ws.Range("J8").Value = 0.1926595011705215
$ws.Range("O8").Value = 0.09959606124235393
$ws.Range("P8").Value = 0.09973689057741834
$ws.Range("Q8").Value = 0.228415196
$ws.Range("R8").Value = 2.055736764
$ws.Range("S8").Value = 0.01918812747750062
$ws.Range("T8").Value = 0.01921525958694431

# Row 9
$ws.Range("D9").Value = "Inflammatory-Mac"
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.11623
$ws.Range("H9").Value = 0.34869
$ws.Range("I9").Value = 0.1926595011705215
$ws.Range("J9").Value = 0.1926595011705215
$ws.Range("M9").Value = 13.379326
$ws.Range("N9").Value = 40.137978
$ws.Range("O9").Value = 0.6780623710957755
$ws.Range("P9").Value = 0.679021154722984
$ws.Range("Q9").Value = 1.55507906098
$ws.Range("R9").Value = 13.99571154882
$ws.Range("S9").Value = 0.1306351581778132
$ws.Range("T9").Value = 0.1308198769531616

# Row 10
$ws.Range("D10").Value = "MuSCs"
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 0.3333333333333333
$ws.Range("G10").Value = 0.11623
$ws.Range("H10").Value = 0.34869
$ws.Range("I10").Value = 0.1926595011705215
$ws.Range("J10").Value = 0.1926595011705215
$ws.Range("M10").Value = 0.08358400000000001
$ws.Range("N10").Value = 0.167168
$ws.Range("O10").Value = 0.004236025433991914
$ws.Range("P10").Value = 0.00282801013027442
$ws.Range("Q10").Value = 0.009714968320000001
$ws.Range("R10").Value = 0.05828980992
$ws.Range("S10").Value = 0.0008161105470585242
$ws.Range("T10").Value = 0.0005448430210038514

# Row 11
$ws.Range("D11").Value = "Resolving-Mac"
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 0.3333333333333333
$ws.Range("G11").Value = 0.11623
$ws.Range("H11").Value = 0.34869
$ws.Range("I11").Value = 0.1926595011705215
$ws.Range("J11").Value = 0.1926595011705215
$ws.Range("M11").Value = 4.270900666666667
$ws.Range("N11").Value = 12.812702
$ws.Range("O11").Value = 0.2164486486654506
$ws.Range("P11").Value = 0.2167547081510057
$ws.Range("Q11").Value = 0.4964067844866667
$ws.Range("R11").Value = 4.46766106038
$ws.Range("S11").Value = 0.04170088868091919
$ws.Range("T11").Value = 0.04175985394873474
